$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "aaaa"
$ws.Range("B11").Value = 1
$ws.Range("D11").Value = 20

$ws.Range("D11").Select()
